$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column P (year 2023) mirroring column O's formatting ---
# Copy number/style formatting from the existing O column cells (rows 3-14)
# down into the new P column, then set the 2023 values.
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

for ($r = 4; $r -le 14; $r++) {
    $ws.Range("O$r").Copy()
    $ws.Range("P$r").PasteSpecial(-4122)
}

$ws.Range("P4").Value = 2023
$ws.Range("P5").Value = 48.2
$ws.Range("P6").Value = 8.6767564891727478
$ws.Range("P7").Value = 12.226605469730881
$ws.Range("P8").Value = 78.520866131691164
$ws.Range("P9").Value = 59.466452648968115
$ws.Range("P10").Value = 26.635270208942913
$ws.Range("P11").Value = 8.166450559693871
$ws.Range("P12").Value = 74.601894583630667
$ws.Range("P13").Value = 99.168063426054971
$ws.Range("P14").Value = 70.956108992253434

# --- Row 14: D14/E14 were blank; now show a right-aligned dash "-" ---
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("D14:E14").HorizontalAlignment = -4152

# --- Row height tweaks ---
$ws.Rows.Item(4).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 15
$ws.Rows.Item(15).RowHeight = 13.5

# --- Clear the lingering selection so the sheet view doesn't point at P8 ---
[void]$ws.Range("A1").Select()
